$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F2").Value = 9.824647717578905
$ws.Range("G2").Value = 9.621358413695905
$ws.Range("H2").Value = 10.02801546946508
$ws.Range("I2").Value = 0.002528524480027724
$ws.Range("J2").Value = 0.002254026368617062
$ws.Range("K2").Value = 0.00286336412815012
$ws.Range("L2").Value = 0.00889024431495853
$ws.Range("M2").Value = 0.008696422842117978
$ws.Range("N2").Value = 0.009091827824325235

$ws.Range("F3").Value = 0.04664409144174743
$ws.Range("G3").Value = 0.04635600292339141
$ws.Range("H3").Value = 0.0469438589313893
$ws.Range("I3").Value = 0.0450843268813746
$ws.Range("J3").Value = 0.04480689623509305
$ws.Range("K3").Value = 0.04537202696875812
$ws.Range("L3").Value = 0.04669460331861312
$ws.Range("M3").Value = 0.04640654441658189
$ws.Range("N3").Value = 0.04699438929953247

$ws.Range("F4").Value = 9.871291809020653
$ws.Range("G4").Value = 9.667714416619296
$ws.Range("H4").Value = 10.07495932839647
$ws.Range("I4").Value = 0.04761285136140232
$ws.Range("J4").Value = 0.04706092260371012
$ws.Range("K4").Value = 0.04823539109690824
$ws.Range("L4").Value = 0.05558484763357165
$ws.Range("M4").Value = 0.05510296725869987
$ws.Range("N4").Value = 0.05608621712385772
